# Format Key Metrics Values on Dashboard.
# Update the Key Metrics row (row 4, "AA") values on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = 4.5
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 5
$ws.Range("S4").Value = 5
